$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1) Header / summary figures that changed value (text labels are unchanged,
#    only the underlying numbers were updated).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 334916      # VALOR MORA total
$ws.Range("F13").Value = 7           # Cant. Periodos

# ---------------------------------------------------------------------------
# 2) The worker table grew from 6 data rows (16-21) to 7 data rows (16-22).
#    Insert a new row at 21 (pushing the former row 21 down to row 22, and
#    every row below - including the signature block rows 26/27 - down by one
#    as well, which is exactly what the target layout needs).
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Insert()

# Give the freshly inserted row the same look/formatting as the data rows
# above it (row 20) before filling in its values.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Rewrite the worker table (rows 16-22) with the updated data: the new
#    worker JOSE JAVIER ALVIS OROZCO (doc 1007390527) with periods 1610-1611,
#    followed by the existing worker NUMAR DEL CRISTO GUERRERO DE AVILA
#    (doc 73121714) with periods 2503, 2505, 2506, 2507 and the new 2508.
# ---------------------------------------------------------------------------
$data = @(
    @(16, "CC", "1007390527", "JOSE JAVIER ALVIS OROZCO",          "1610", 27578, 689455),
    @(17, "CC", "1007390527", "JOSE JAVIER ALVIS OROZCO",          "1611", 27578, 689455),
    @(18, "CC", "73121714",   "NUMAR DEL CRISTO GUERRERO DE AVILA", "2503", 52000, 1423500),
    @(19, "CC", "73121714",   "NUMAR DEL CRISTO GUERRERO DE AVILA", "2505", 56940, 1423500),
    @(20, "CC", "73121714",   "NUMAR DEL CRISTO GUERRERO DE AVILA", "2506", 56940, 1423500),
    @(21, "CC", "73121714",   "NUMAR DEL CRISTO GUERRERO DE AVILA", "2507", 56940, 1423500),
    @(22, "CC", "73121714",   "NUMAR DEL CRISTO GUERRERO DE AVILA", "2508", 56940, 1423500)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}

Write-Host "Edit applied"
